# "Fruta / hortaliza, semanal" weekly data refresh.
# A new weekly record for "Arveja Verde" (Vega Modelo de Temuco) needs to be
# inserted as row 62, pushing the existing rows 62-77 down to 63-78 and
# extending the used range from A1:R77 to A1:R78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62 (shifts rows 62..77 down to 63..78).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value = 10
$ws.Cells.Item(62, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(62, 3).Value = "La Araucanía"
$ws.Cells.Item(62, 4).Value = 44551
$ws.Cells.Item(62, 5).Value = 9
$ws.Cells.Item(62, 6).Value = 100112022
$ws.Cells.Item(62, 7).Value = "Arveja Verde"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 155
$ws.Cells.Item(62, 11).Value = 15000
$ws.Cells.Item(62, 12).Value = 15000
$ws.Cells.Item(62, 13).Value = 15000
$ws.Cells.Item(62, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(62, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(62, 16).Value = 600
$ws.Cells.Item(62, 17).Value = 25
$ws.Cells.Item(62, 18).Value = "Hortaliza"
